$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the Quantidade values: they were stored as text ("420,00" / "610,00")
# and must become plain numbers (420 / 610).
$ws.Range("E2").Value = 420
$ws.Range("E3").Value = 610

# Remove the conditional formatting (duplicate values rule) on column A.
$ws.Cells.FormatConditions.Delete()

# Re-apply the AutoFilter so its range matches the actual used data range
# instead of the whole column.
$ws.AutoFilterMode = $false
$ws.Range("A1:E3").AutoFilter()

$wb.Save()
